# Apply the Booking Report corrections:
#  - fix the "confirmedd" typo in the Status column
#  - delete the trailing duplicate rows (rows 8-10)
#  - renumber the Booking ID column continuously (1..6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in row 2's Status cell.
$ws.Range("F2").Value = "confirmed"

# Remove the three trailing junk/duplicate rows (originally rows 8, 9, 10).
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()

# Renumber the Booking ID column continuously for the remaining data rows.
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
